$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Rename the stray "asdasd" text (shared string, currently only referenced by
# A50) to "cada 1 km" in place, then relocate it to the new summary table
# header at M1, and drop the now-empty A50 cell/row.
$ws.Range("A50").Value = "cada 1 km"
$etiqueta = $ws.Range("A50").Value2
$ws.Range("M1").Value = $etiqueta
[void]$ws.Rows.Item(50).Delete()

# New "cada 1 km" summary table header (K2:N2) mirrors the A2:D2 headers.
$ws.Range("K2").Value = "de"
$ws.Range("L2").Value = "a"
$ws.Range("M2").Value = "corte"
$ws.Range("N2").Value = "relleno"

# Bucket boundaries (every 1000 units) and per-bucket totals of corte (C) /
# relleno (D), 10 source rows per bucket.
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 1000
$ws.Range("M3").Formula = "=SUM(C3:C12)"
$ws.Range("N3").Formula = "=SUM(D3:D12)"

$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 2000
$ws.Range("M4").Formula = "=SUM(C13:C22)"
$ws.Range("N4").Formula = "=SUM(D13:D22)"

$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 3000
$ws.Range("M5").Formula = "=SUM(C23:C32)"
$ws.Range("N5").Formula = "=SUM(D23:D32)"

$ws.Range("K6").Value = 3000
$ws.Range("L6").Value = 4000
$ws.Range("M6").Formula = "=SUM(C33:C42)"
$ws.Range("N6").Formula = "=SUM(D33:D42)"

$ws.Range("K7").Value = 4000
$ws.Range("L7").Value = 4300
$ws.Range("M7").Formula = "=SUM(C43:C45)"
$ws.Range("N7").Formula = "=SUM(D43:D45)"

# Grand total of the bucketed "corte" column.
$ws.Range("M8").Formula = "=SUM(M3:M7)"

# Leave the selection where the new table was built.
[void]$ws.Range("L8").Select()
